# Generate Report for handoff
# Updates the "Latest Handoff Datetime" for the row whose source file is
# 6fdc26cb-b1cd-4a2f-9353-7cf3f4da1435.md (status "Ready for handoff"),
# on both locale report sheets, reflecting a fresh handoff just performed.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-18 06:20:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-18 06:20:10"
